# Generate Report for Handback
# The last handoff file (b7c20352-...) failed the handback transform because
# the handback file name (m2u5ltcf.vgj) did not match the handoff file name.
# Record the failure on the Overview sheet and on each language sheet, and
# widen the "Error Detail" column so the message is readable.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$failedStatus = "Handback transform failed"

# Overview sheet: row 7 is the b7c20352-61ee-4299-98e8-ad8b9d197f17.md file;
# both the zh-cn (E) and de-de (F) status columns move to the failure state.
$overview.Range("E7").Value = $failedStatus
$overview.Range("F7").Value = $failedStatus

# zh-cn sheet: row 7 Status + Error Detail
$zhcn.Range("C7").Value = $failedStatus
$zhcn.Range("P7").Value = "Handback file name: m2u5ltcf.vgj is different with handoff file name: b7c20352-61ee-4299-98e8-ad8b9d197f17.6d9f45156d46bfdb16e11ba4b4676eb186f7e6de.zh-cn."

# de-de sheet: row 7 Status + Error Detail
$dede.Range("C7").Value = $failedStatus
$dede.Range("P7").Value = "Handback file name: m2u5ltcf.vgj is different with handoff file name: b7c20352-61ee-4299-98e8-ad8b9d197f17.6d9f45156d46bfdb16e11ba4b4676eb186f7e6de.de-de."

# Widen column P ("Error Detail") on both language sheets so the new message is visible.
# (39.16 is the COM ColumnWidth input that round-trips to a stored OOXML column
# width of exactly 40 through Excel's pixel-quantized width model.)
$zhcn.Columns.Item(16).ColumnWidth = 39.16
$dede.Columns.Item(16).ColumnWidth = 39.16
